$wb = $excel.ActiveWorkbook

# Update the "Metadata" sheet
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/room-and-board-service"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# Update the "Include from Room And Board S" sheet
$include = $wb.Worksheets.Item("Include from Room And Board S")
$include.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/room-and-board-service"
